$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crediti a inizi 2023")

# --- Insert the new "Name" / "Value" header row above the old row 4 ---
# (pushes "type / TradeReceivableCredits" from row4->row5, "value column / date" from row5->row6, etc.)
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Style = "Normal"
$ws.Range("A4").Value = "Name"
$ws.Range("B4").Value = "Value"
$ws.Range("A4").Font.Size = 12
$ws.Rows.Item(4).RowHeight = 15.5

# --- Insert the new "opposite type" / "Cash" data row below the (shifted) "value column" row ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Style = "Normal"
$ws.Range("A7").Value = "opposite type"
$ws.Range("B7").Value = "Cash"

# --- The existing "CO__crediti2022" table has shifted down by two rows; move it to match ---
$creditiTable = $ws.ListObjects.Item("CO__crediti2022")
$creditiTable.Resize($ws.Range("A16:E18"))

# --- Turn the new A4:B7 block into its own table ("Table5") ---
$settingsTable = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A4:B7"), [System.Type]::Missing, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$settingsTable.Name = "Table5"

# --- Update the sheet selection ---
$ws.Range("A7").Select()
